# Generate Report for Handoff
# Update the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect that the file is now ready for handoff (instead of "handed back"),
# with fresh handoff file names / timestamps and a new error detail message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3b266b929e1d6479e408150d78f672dc4feca97/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6cd6fe211439830c9d40b8d7ffbd034071b79c54/e2e/b.md."

# ---------------------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-07 10:54:16"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 10:54:02"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 10:54:16"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 40
